$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks numeric (e.g. "-0.014") as genuine text
# rather than letting Excel auto-convert it to a number. We do this by
# entering it as a text formula (string concatenation always yields text),
# then flattening the formula to a static value with a values-only paste.
# This keeps the cell's type as a shared string (t="s") without leaving
# any NumberFormat/style residue behind, unlike a quote-prefix or "@"
# number-format approach would.
function Set-TextValue($cell, $text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# Row 2 (C/A Lag): -0.549* / -0.023  ->  -0.546* / -0.014
$ws.Range("B2").Value = "-0.546*"
Set-TextValue $ws.Range("C2") "-0.014"

# Row 3 (LF Lag): 0.906 / 0.349  ->  1.055** / 0.922***
$ws.Range("B3").Value = "1.055**"
$ws.Range("C3").Value = "0.922***"

# Row 4 previously held "Constant" / -0.217 / -0.834***.
# It becomes the "r2" row, with numeric (not text) values, taken from the
# old row 5 ("r2_adj") which is being removed.
$ws.Range("A4").Value = "r2"
$ws.Range("B4").Value = 0.2032637665075804
$ws.Range("C4").Value = 0.3626403328286422

# Old row 5 ("r2_adj" / 0.13 / 0.01) is deleted entirely, shifting nothing
# else (it was the last row).
$ws.Rows.Item(5).Delete()
